# Adds the "route 109" and "seashore house" trainer blocks to the trainers
# sheet, moving the trailing "END" marker down to make room for them.
#
# Before: ... TRAINER_BRAWLY_1 block ends at row 398, then a gap, then
#         "END" sits alone at row 404 (dimension A1:G404).
# After:  "END" moves to row 400, and rows 402-448 hold the new
#         "# route 109" and "# seashore house" trainer entries
#         (dimension A1:G448).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- END marker moves from A404 up to A400 ---
$ws.Range("A400").Value = "END"

# --- # route 109 ---
$ws.Range("A402").Value = "# route 109"

$ws.Range("A403").Value = "species"
$ws.Range("B403").Value = "lvl"
$ws.Range("C403").Value = "iv"
$ws.Range("D403").Value = "heldItem"
$ws.Range("E403").Value = "moves"
$ws.Range("F403").Value = "ability"
$ws.Range("G403").Value = "shiny"

$ws.Range("A404").Value = "Machop"
$ws.Range("B404").Value = 24

$ws.Range("A405").Value = "TRAINER_HUEY"
$ws.Range("B405").Value = 24

$ws.Range("A407").Value = "Poliwhirl"

$ws.Range("A408").Value = "species"
$ws.Range("B408").Value = "lvl"
$ws.Range("C408").Value = "iv"
$ws.Range("D408").Value = "heldItem"
$ws.Range("E408").Value = "moves"
$ws.Range("F408").Value = "ability"
$ws.Range("G408").Value = "shiny"

$ws.Range("A409").Value = "Azurill"
$ws.Range("B409").Value = 24

$ws.Range("A410").Value = "TRAINER_HAILEY"
$ws.Range("B410").Value = 25

$ws.Range("A412").Value = "Marill"

$ws.Range("A413").Value = "species"
$ws.Range("B413").Value = "lvl"
$ws.Range("C413").Value = "iv"
$ws.Range("D413").Value = "heldItem"
$ws.Range("E413").Value = "moves"
$ws.Range("F413").Value = "ability"
$ws.Range("G413").Value = "shiny"

$ws.Range("A414").Value = "TRAINER_EDMOND"
$ws.Range("B414").Value = 25

$ws.Range("A415").Value = "Seel"
$ws.Range("B415").Value = 25

$ws.Range("A417").Value = "Lombre"

$ws.Range("A418").Value = "species"
$ws.Range("B418").Value = "lvl"
$ws.Range("C418").Value = "iv"
$ws.Range("D418").Value = "heldItem"
$ws.Range("E418").Value = "moves"
$ws.Range("F418").Value = "ability"
$ws.Range("G418").Value = "shiny"

$ws.Range("A419").Value = "Mudkip"
$ws.Range("B419").Value = 24

$ws.Range("A420").Value = "Wooper"
$ws.Range("B420").Value = 25

$ws.Range("A422").Value = "TRAINER_RICKY_1"

$ws.Range("A423").Value = "species"
$ws.Range("B423").Value = "lvl"
$ws.Range("C423").Value = "iv"
$ws.Range("D423").Value = "heldItem"
$ws.Range("E423").Value = "moves"
$ws.Range("F423").Value = "ability"
$ws.Range("G423").Value = "shiny"

$ws.Range("A424").Value = "TRAINER_LOLA_1"
$ws.Range("B424").Value = 25
$ws.Range("F424").Value = "Masquerain"

$ws.Range("A426").Value = "Flutter"

$ws.Range("A427").Value = "species"
$ws.Range("B427").Value = "lvl"
$ws.Range("C427").Value = "iv"
$ws.Range("D427").Value = "heldItem"
$ws.Range("E427").Value = "moves"
$ws.Range("F427").Value = "ability"
$ws.Range("G427").Value = "shiny"

$ws.Range("A428").Value = "TRAINER_CHANDLER"
$ws.Range("B428").Value = 24

$ws.Range("A429").Value = "Sandshrew"
$ws.Range("B429").Value = 25

$ws.Range("A431").Value = "Wingull"

# --- # seashore house ---
$ws.Range("A433").Value = "Lombre"

$ws.Range("A434").Value = "species"
$ws.Range("B434").Value = "lvl"
$ws.Range("C434").Value = "iv"
$ws.Range("D434").Value = "heldItem"
$ws.Range("E434").Value = "moves"
$ws.Range("F434").Value = "ability"
$ws.Range("G434").Value = "shiny"

$ws.Range("A435").Value = "Marshtomp"
$ws.Range("B435").Value = 26

$ws.Range("A436").Value = "# seashore house"
$ws.Range("B436").Value = 28

$ws.Range("A438").Value = "Machoke"

$ws.Range("A439").Value = "species"
$ws.Range("B439").Value = "lvl"
$ws.Range("C439").Value = "iv"
$ws.Range("D439").Value = "heldItem"
$ws.Range("E439").Value = "moves"
$ws.Range("F439").Value = "ability"
$ws.Range("G439").Value = "shiny"

$ws.Range("A440").Value = "TRAINER_JOHANNA"
$ws.Range("B440").Value = 26

$ws.Range("A441").Value = "Linoone"
$ws.Range("B441").Value = 27

$ws.Range("A443").Value = "Raichu"

$ws.Range("A444").Value = "species"
$ws.Range("B444").Value = "lvl"
$ws.Range("C444").Value = "iv"
$ws.Range("D444").Value = "heldItem"
$ws.Range("E444").Value = "moves"
$ws.Range("F444").Value = "ability"
$ws.Range("G444").Value = "shiny"

$ws.Range("A445").Value = "TRAINER_SIMON"
$ws.Range("B445").Value = 25

$ws.Range("A446").Value = "Sentret"
$ws.Range("B446").Value = 26

$ws.Range("A448").Value = "Psyduck"

# --- restore the selection/cursor position recorded in the saved file ---
$ws.Range("B451").Select()
